$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pull_subscription")

# Update E5 value from 30 to 120
$ws.Range("E5").Value = 120

# Populate row 6 with new data
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "pull-sub-10"
$ws.Range("C6").Value = "pull-sub-10"
$ws.Range("D6").Value = "topic-10"
$ws.Range("E6").Value = 60
$ws.Range("F6").Value = "2400s"
$ws.Range("G6").Value = 5

# Remove row 7 entirely (it was empty placeholder row)
$ws.Rows.Item(7).Delete()
